# Update column E (Custom Betweenness (R=actual)) values on Sheet1
# to reflect corrected figures found after testing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 78507.97
    3  = 31460.08
    4  = 27657.51
    5  = 116789.4
    6  = 71157.02
    8  = 12877.46
    9  = 79904.2
    10 = 13325.17
    11 = 70035.75999999999
    13 = 27127.73
    14 = 6620.4
    16 = 13025.25
    17 = 2248.97
    18 = 21695.74
    19 = 4035.65
    20 = 1957.3
    21 = 16817.6
    22 = 19904
    23 = 14026.26
    24 = 12130.22
    25 = 4847.03
    26 = 16521.27
    27 = 5328.96
    30 = 5343.89
    31 = 5047.9
    32 = 9225.91
    34 = 10831.5
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 5).Value = $updates[$row]
}

$wb.Save()
